$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 90
$ws.Range("C2").Value = 'face/face002.jpg'
$ws.Range("D2").Value = 'töten'
$ws.Range("E2").Value = 'face'

$ws.Range("B3").Value = 52
$ws.Range("C3").Value = 'house/house008.jpg'
$ws.Range("D3").Value = 'schicken'
$ws.Range("E3").Value = 'house'

$ws.Range("B4").Value = 124
$ws.Range("C4").Value = 'house/house004.jpg'
$ws.Range("D4").Value = 'tauschen'
$ws.Range("E4").Value = 'house'

$ws.Range("B5").Value = 20
$ws.Range("C5").Value = 'house/house014.jpg'
$ws.Range("D5").Value = 'stärken'
$ws.Range("E5").Value = 'house'

$ws.Range("B6").Value = 29
$ws.Range("C6").Value = 'face/face012.jpg'
$ws.Range("D6").Value = 'währen'
$ws.Range("E6").Value = 'face'

$ws.Range("B7").Value = 117
$ws.Range("C7").Value = 'face/face024.jpg'
$ws.Range("D7").Value = 'fliehen'
$ws.Range("E7").Value = 'face'

$ws.Range("B8").Value = 7
$ws.Range("C8").Value = 'house/house013.jpg'
$ws.Range("D8").Value = 'schätzen'
$ws.Range("E8").Value = 'house'

$ws.Range("B9").Value = 25
$ws.Range("C9").Value = 'house/house022.jpg'
$ws.Range("D9").Value = 'kehren'
$ws.Range("E9").Value = 'house'

$ws.Range("B10").Value = 71
$ws.Range("C10").Value = 'house/house024.jpg'
$ws.Range("D10").Value = 'antun'
$ws.Range("E10").Value = 'house'

$ws.Range("B11").Value = 59
$ws.Range("C11").Value = 'face/face014.jpg'
$ws.Range("D11").Value = 'fühlen'
$ws.Range("E11").Value = 'face'

$ws.Range("B12").Value = 15
$ws.Range("C12").Value = 'house/house025.jpg'
$ws.Range("D12").Value = 'opfern'
$ws.Range("E12").Value = 'house'

$ws.Range("B13").Value = 32
$ws.Range("C13").Value = 'house/house015.jpg'
$ws.Range("D13").Value = 'posten'
$ws.Range("E13").Value = 'house'

$ws.Range("B14").Value = 94
$ws.Range("C14").Value = 'house/house010.jpg'
$ws.Range("D14").Value = 'fliegen'
$ws.Range("E14").Value = 'house'

$ws.Range("B15").Value = 96
$ws.Range("C15").Value = 'face/face028.jpg'
$ws.Range("D15").Value = 'hoffen'
$ws.Range("E15").Value = 'face'

$ws.Range("B16").Value = 31
$ws.Range("C16").Value = 'house/house017.jpg'
$ws.Range("D16").Value = 'klappen'
$ws.Range("E16").Value = 'house'

$ws.Range("B17").Value = 113
$ws.Range("C17").Value = 'house/house030.jpg'
$ws.Range("D17").Value = 'husten'
$ws.Range("E17").Value = 'house'

$ws.Range("B18").Value = 103
$ws.Range("C18").Value = 'house/house026.jpg'
$ws.Range("D18").Value = 'bleiben'
$ws.Range("E18").Value = 'house'

$ws.Range("B19").Value = 75
$ws.Range("C19").Value = 'face/face023.jpg'
$ws.Range("D19").Value = 'sondern'
$ws.Range("E19").Value = 'face'

$ws.Range("B20").Value = 85
$ws.Range("C20").Value = 'house/house003.jpg'
$ws.Range("D20").Value = 'segeln'
$ws.Range("E20").Value = 'house'

$ws.Range("B21").Value = 92
$ws.Range("C21").Value = 'face/face003.jpg'
$ws.Range("D21").Value = 'rasen'
$ws.Range("E21").Value = 'face'

$ws.Range("B22").Value = 21
$ws.Range("C22").Value = 'face/face021.jpg'
$ws.Range("D22").Value = 'nehmen'
$ws.Range("E22").Value = 'face'

$ws.Range("B23").Value = 58
$ws.Range("C23").Value = 'house/house016.jpg'
$ws.Range("D23").Value = 'schenken'
$ws.Range("E23").Value = 'house'

$ws.Range("B24").Value = 48
$ws.Range("C24").Value = 'face/face029.jpg'
$ws.Range("D24").Value = 'liefern'
$ws.Range("E24").Value = 'face'

$ws.Range("B25").Value = 78
$ws.Range("C25").Value = 'house/house027.jpg'
$ws.Range("D25").Value = 'formen'
$ws.Range("E25").Value = 'house'

$ws.Range("B26").Value = 60
$ws.Range("C26").Value = 'face/face013.jpg'
$ws.Range("D26").Value = 'schmecken'
$ws.Range("E26").Value = 'face'

$ws.Range("B27").Value = 107
$ws.Range("C27").Value = 'face/face019.jpg'
$ws.Range("D27").Value = 'enden'
$ws.Range("E27").Value = 'face'

$ws.Range("B28").Value = 19
$ws.Range("C28").Value = 'face/face025.jpg'
$ws.Range("D28").Value = 'sieben'
$ws.Range("E28").Value = 'face'

$ws.Range("B29").Value = 36
$ws.Range("C29").Value = 'face/face020.jpg'
$ws.Range("D29").Value = 'haken'
$ws.Range("E29").Value = 'face'

$ws.Range("B30").Value = 111
$ws.Range("C30").Value = 'face/face007.jpg'
$ws.Range("D30").Value = 'hauen'
$ws.Range("E30").Value = 'face'

$ws.Range("B31").Value = 100
$ws.Range("C31").Value = 'face/face027.jpg'
$ws.Range("D31").Value = 'kaufen'
$ws.Range("E31").Value = 'face'

$ws.Range("B32").Value = 11
$ws.Range("C32").Value = 'house/house002.jpg'
$ws.Range("D32").Value = 'raten'
$ws.Range("E32").Value = 'house'

$ws.Range("B33").Value = 45
$ws.Range("C33").Value = 'face/face005.jpg'
$ws.Range("D33").Value = 'spielen'
$ws.Range("E33").Value = 'face'
